$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "LC69055-76"
$ws.Range("B2").Value = "C:/Users/Asus/Desktop/тест пнг\LC69055-76\5_image_LC69055-76.jpg"
$ws.Range("C2").Value = 977

$ws.Range("A3").Value = "LC69802-01"
$ws.Range("B3").Value = "C:/Users/Asus/Desktop/тест пнг\LC69802-01\5_image_LC69802-01.jpg"
$ws.Range("C3").Value = 1154
